# Apply updated crypto price/volume data to Sheet1 (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.156.45"
$ws.Range("E2").Value = "  -4.27%  "

$ws.Range("D3").Value = "1.652.28"
$ws.Range("E3").Value = "  -3.77%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'215.61"
$ws.Range("E5").Value = "  -4.40%  "

$ws.Range("D6").Value = "'0.5130"
$ws.Range("E6").Value = "  -3.25%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").Value = "'0.2594"
$ws.Range("E8").Value = "  -2.77%  "

$ws.Range("D9").Value = "'0.06440"
$ws.Range("E9").Value = "  -3.89%  "

$ws.Range("D10").Value = "'19.93"
$ws.Range("E10").Value = "  -4.86%  "

$ws.Range("D11").Value = "'0.07780"
$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("D12").Value = "1.663.23"
$ws.Range("E12").Value = "  -3.61%  "

$ws.Range("E13").Value = "  -4.86%  "

$ws.Range("D14").Value = "1.878.66"
$ws.Range("E14").Value = "  -3.79%  "

$ws.Range("D15").Value = "'0.5521"
$ws.Range("E15").Value = "  -5.92%  "

$ws.Range("D16").Value = "0.0₅8001"
$ws.Range("E16").Value = "  -2.69%  "

$ws.Range("D18").Value = "26.164.89"
$ws.Range("E18").Value = "  -4.56%  "

$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").Value = "'210.72"
$ws.Range("E20").Value = "  -5.54%  "

$ws.Range("D21").Value = "'4.400"
$ws.Range("E21").Value = "  -5.67%  "

$ws.Range("E22").Value = "  -4.34%  "

$ws.Range("D23").Value = "'6.048"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").Value = "'143.45"
$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("D26").Value = "'1.749"
$ws.Range("E26").Value = "  +3.18%  "

$ws.Range("D27").Value = "'0.1176"
$ws.Range("E27").Value = "  -2.82%  "

$ws.Range("D28").Value = "'6.970"
$ws.Range("E28").Value = "  -3.97%  "

$ws.Range("E29").Value = "  -2.83%  "

$ws.Range("D30").Value = "'0.05091"
$ws.Range("E30").Value = "  -5.12%  "

$ws.Range("E31").Value = "  -4.12%  "

$ws.Range("D32").Value = "'3.355"

$ws.Range("D33").Value = "'3.221"
$ws.Range("E33").Value = "  -6.24%  "

$ws.Range("E34").Value = "  -4.31%  "

$ws.Range("D35").Value = "'2.738"
$ws.Range("E35").Value = "  -4.62%  "

$ws.Range("D36").Value = "'0.9238"
$ws.Range("E36").Value = "  -3.63%  "

$ws.Range("D37").Value = "'2.350"
$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("D38").Value = "'0.5714"
$ws.Range("E38").Value = "  -2.67%  "

$ws.Range("D39").Value = "1.161.54"
$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("D40").Value = "'0.01586"
$ws.Range("E40").Value = "  -3.80%  "

$ws.Range("D41").Value = "'2.560"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").Value = "'1.005"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("D43").Value = "'5.654"
$ws.Range("E43").Value = "  -2.59%  "

$ws.Range("D44").Value = "'0.8231"
$ws.Range("E44").Value = "  -2.29%  "

$ws.Range("D45").Value = "'100.09"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("D46").Value = "1.789.49"
$ws.Range("E46").Value = "  -3.72%  "

$ws.Range("D47").Value = "0.0₈116"
$ws.Range("E47").Value = "  +5.15%  "

$ws.Range("D48").Value = "'0.4550"
$ws.Range("E48").Value = "  -0.94%  "

$ws.Range("E49").Value = "  -4.02%  "

$ws.Range("D50").Value = "'1.005"

$ws.Range("D51").Value = "'7.857"
$ws.Range("E51").Value = "  -3.40%  "
